$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 424, shifting existing rows 424.. down by one.
$ws.Rows("424:424").Insert()

# Populate the newly inserted row 424 with the new data record.
$ws.Cells.Item(424, 1).Value = 6
$ws.Cells.Item(424, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(424, 3).Value = "Metropolitana"
$ws.Cells.Item(424, 4).Value = 44889
$ws.Cells.Item(424, 5).Value = 13
$ws.Cells.Item(424, 6).Value = 100112030
$ws.Cells.Item(424, 7).Value = "Poroto granado"
$ws.Cells.Item(424, 8).Value = "Sin especificar"
$ws.Cells.Item(424, 9).Value = "Primera"
$ws.Cells.Item(424, 10).Value = 400
$ws.Cells.Item(424, 11).Value = 40000
$ws.Cells.Item(424, 12).Value = 42000
$ws.Cells.Item(424, 13).Value = 40850
$ws.Cells.Item(424, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(424, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(424, 16).Value = 2723
$ws.Cells.Item(424, 17).Value = 15
$ws.Cells.Item(424, 18).Value = "Hortaliza"
